$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style (bold, border, centered) from H1 to the new I1:J1 header cells
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# New header labels
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data columns I0 (col I) and IF (col J) for rows 2-74
$data = @(
  @(7,8), @(7,8), @(2,2), @(9,9), @(7,8), @(6,6), @(7,8), @(6,6), @(1,2), @(1,3),
  @(8,8), @(9,9), @(9,9), @(7,7), @(6,7), @(11,11), @(9,9), @(9,9), @(8,8), @(8,8),
  @(8,8), @(8,8), @(8,8), @(9,10), @(7,8), @(9,9), @(9,9), @(8,9), @(9,9), @(8,9),
  @(9,9), @(9,9), @(9,10), @(9,9), @(8,8), @(7,7), @(8,8), @(7,8), @(7,8), @(7,8),
  @(8,8), @(6,7), @(10,12), @(8,8), @(7,8), @(7,8), @(8,9), @(9,9), @(6,6), @(8,8),
  @(7,7), @(7,7), @(6,7), @(7,8), @(6,8), @(7,8), @(7,7), @(9,9), @(8,8), @(7,8),
  @(9,9), @(9,9), @(8,9), @(9,9), @(6,6), @(8,8), @(9,9), @(7,7), @(6,7), @(6,6),
  @(4,4), @(4,4), @(4,4)
)

for ($idx = 0; $idx -lt $data.Length; $idx++) {
    $row = $idx + 2
    $pair = $data[$idx]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
